$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.813.67'
$ws.Range("E2").Value = '  -2.25%  '

$ws.Range("D3").Value = '3.851.13'
$ws.Range("E3").Value = '  -2.76%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.10'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.667'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.31%  '

$ws.Range("E8").Value = '  +0.27%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.747'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.31%  '

$ws.Range("E10").Value = '  +4.22%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.04'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000320'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.81%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.85%  '

$ws.Range("D14").Value = '4.483.39'
$ws.Range("E14").Value = '  -2.28%  '

$ws.Range("D15").Value = '3.861.15'
$ws.Range("E15").Value = '  -3.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.59'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.72%  '

$ws.Range("E18").Value = '  -6.11%  '

$ws.Range("E19").Value = '  -2.13%  '

$ws.Range("D20").Value = '70.833.06'
$ws.Range("E20").Value = '  -1.96%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '432.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '93.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.82%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.35%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.82%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.09'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.40%  '

$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("E29").Value = '  -4.62%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.88'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.46'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.88%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '48.85'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.124'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.64%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '69.22'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.75%  '

$ws.Range("D36").Value = '0.0₃0967'
$ws.Range("E36").Value = '  +14.07%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '618.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.417'
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.26%  '

$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.30%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.142'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.98%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +31.64%  '

$ws.Range("B43").Value = 'ThetaToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.25'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.18%  '

$ws.Range("E44").Value = '  -3.66%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.66'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.41%  '

$ws.Range("E47").Value = '  -2.99%  '

$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.80'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -16.36%  '

$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.16%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.822.35'
$ws.Range("E50").Value = '  +2.54%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000271'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.35%  '

